# Fruta / hortaliza, semanal
#
# A new weekly price record needs to be inserted into the "Sandia"
# (Femacal de La Calera) sheet. The new record belongs chronologically
# before the row that is currently at row 650, so we insert a brand new
# row there (pushing everything from 650 downwards by one row) and then
# populate it with the new observation's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 650; rows 650..717 shift down to 651..718
$ws.Rows("650:650").Insert()

$ws.Range("A650").Value = 3
$ws.Range("B650").Value = "Femacal de La Calera"
$ws.Range("C650").Value = "Coquimbo"
$ws.Range("D650").Value = 45194
$ws.Range("E650").Value = 5
$ws.Range("F650").Value = 100112028
$ws.Range("G650").Value = "Sandia"
$ws.Range("H650").Value = "Sin especificar"
$ws.Range("I650").Value = "Primera"
$ws.Range("J650").Value = 160
$ws.Range("K650").Value = 1000
$ws.Range("L650").Value = 1000
$ws.Range("M650").Value = 1000
$ws.Range("N650").Value = "$/kilo (volumen en unidades)"
$ws.Range("O650").Value = "Perú"
$ws.Range("P650").Value = 1000
$ws.Range("Q650").Value = 1
$ws.Range("R650").Value = "Hortaliza"

# Match the date-number formatting used by the other "Fecha" cells in
# column D (style carried over automatically by Insert, but make sure).
$ws.Range("D650").NumberFormat = $ws.Range("D651").NumberFormat
